$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 51
$ws.Range("H51").Value = 2961.577
$ws.Range("I51").Value = 2999.9473
$ws.Range("K51").Value = 2999.9473
$ws.Range("M51").Value = -2515.9473
# Row 60
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
# Row 98
$ws.Range("H98").Value = 1363.7609
$ws.Range("I98").Value = 1430.4634
$ws.Range("K98").Value = 1430.4634
$ws.Range("M98").Value = 67.53659999999991
# Row 122
$ws.Range("H122").Value = 1363.7609
$ws.Range("I122").Value = 1430.4634
$ws.Range("K122").Value = 4291.3902
$ws.Range("M122").Value = -1841.3902
# Row 132
$ws.Range("H132").Value = 4203.906
$ws.Range("I132").Value = 4618.829
$ws.Range("J132").Value = 2786.25
$ws.Range("K132").Value = 13856.487
$ws.Range("L132").Value = 8358.75
$ws.Range("M132").Value = -11326.487
$ws.Range("N132").Value = -13418.75
# Row 136
$ws.Range("H136").Value = 50390
$ws.Range("J136").Value = 50390
$ws.Range("L136").Value = 50390
$ws.Range("N136").Value = -60590

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1892
$ws.Range("I110").Value = 1126.5
$ws.Range("K110").Value = 1126.5
$ws.Range("M110").Value = 918.5
# Row 122
$ws.Range("H122").Value = 5540
$ws.Range("I122").Value = 4691.1665
$ws.Range("J122").Value = 7449.875
$ws.Range("K122").Value = 14073.4995
$ws.Range("L122").Value = 22349.625
$ws.Range("M122").Value = -11623.4995
$ws.Range("N122").Value = -27249.625
# Row 128
$ws.Range("H128").Value = 99884.5
$ws.Range("J128").Value = 99884.5
$ws.Range("L128").Value = 99884.5
$ws.Range("N128").Value = -109844.5
# Row 132
$ws.Range("H132").Value = 2802.9512
$ws.Range("I132").Value = 2009.3334
$ws.Range("K132").Value = 6028.0002
$ws.Range("M132").Value = -3498.0002

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3055.2
$ws.Range("I105").Value = 2755.9
$ws.Range("K105").Value = 2755.9
$ws.Range("M105").Value = -1008.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2741.7646
$ws.Range("I16").Value = 2601.1538
$ws.Range("K16").Value = 2601.1538
$ws.Range("M16").Value = -2314.1538
# Row 75
$ws.Range("H75").Value = 90424.42999999999
$ws.Range("J75").Value = 102794.6
$ws.Range("L75").Value = 102794.6
$ws.Range("N75").Value = -104790.6
# Row 78
$ws.Range("H78").Value = 90424.42999999999
$ws.Range("J78").Value = 102794.6
$ws.Range("L78").Value = 308383.8
$ws.Range("N78").Value = -318367.8
# Row 113
$ws.Range("H113").Value = 2741.7646
$ws.Range("I113").Value = 2601.1538
$ws.Range("K113").Value = 2601.1538
$ws.Range("M113").Value = -431.1538
# Row 122
$ws.Range("H122").Value = 2585.5625
$ws.Range("I122").Value = 2388.182
$ws.Range("K122").Value = 7164.545999999999
$ws.Range("M122").Value = -4714.545999999999
# Row 132
$ws.Range("H132").Value = 2466.5454
$ws.Range("J132").Value = 2608.3333
$ws.Range("L132").Value = 7824.999899999999
$ws.Range("N132").Value = -12884.9999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 1268.4445
$ws.Range("I9").Value = 103.2
$ws.Range("K9").Value = 309.6
$ws.Range("M9").Value = -85.60000000000002
# Row 15
$ws.Range("H15").Value = 1000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 3000
$ws.Range("N15").Value = -3280
$ws.Range("M15").ClearContents()
# Row 20
$ws.Range("H20").Value = 1699.8334
$ws.Range("J20").Value = 1350
$ws.Range("L20").Value = 4050
$ws.Range("N20").Value = -4504
# Row 21
$ws.Range("H21").Value = 250
$ws.Range("I21").Value = 250
$ws.Range("K21").Value = 750
$ws.Range("M21").Value = -577
# Row 40
$ws.Range("H40").Value = 2112.4
$ws.Range("I40").Value = 61.25
$ws.Range("J40").Value = 3479.8333
$ws.Range("K40").Value = 245
$ws.Range("L40").Value = 13919.3332
$ws.Range("M40").Value = -176
$ws.Range("N40").Value = -14057.3332
# Row 113
$ws.Range("H113").Value = 1536.6
$ws.Range("I113").Value = 1828.1666
$ws.Range("K113").Value = 5484.4998
$ws.Range("M113").Value = -3314.4998
# Row 128
$ws.Range("H128").Value = 88087
$ws.Range("I128").Value = 88087
$ws.Range("K128").Value = 264261
$ws.Range("M128").Value = -259281

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 33598
$ws.Range("J57").Value = 44996.668
$ws.Range("L57").Value = 44996.668
$ws.Range("N57").Value = -46636.668
# Row 122
$ws.Range("H122").Value = 2637.2856
$ws.Range("I122").Value = 1666
$ws.Range("J122").Value = 3365.75
$ws.Range("K122").Value = 4998
$ws.Range("L122").Value = 10097.25
$ws.Range("M122").Value = -2548
$ws.Range("N122").Value = -14997.25
# Row 126
$ws.Range("H126").Value = 2362.7273
$ws.Range("J126").Value = 2799
$ws.Range("L126").Value = 8397
$ws.Range("N126").Value = -13337
# Row 132
$ws.Range("H132").Value = 2095.6086
$ws.Range("I132").Value = 1681.4706
$ws.Range("J132").Value = 3269
$ws.Range("K132").Value = 5044.4118
$ws.Range("L132").Value = 9807
$ws.Range("M132").Value = -2514.4118
$ws.Range("N132").Value = -14867

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2508.0908
$ws.Range("I7").Value = 2425.75
$ws.Range("K7").Value = 2425.75
$ws.Range("M7").Value = -2313.75
# Row 126
$ws.Range("H126").Value = 2508.0908
$ws.Range("I126").Value = 2425.75
$ws.Range("K126").Value = 7277.25
$ws.Range("M126").Value = -4807.25
# Row 132
$ws.Range("H132").Value = 3318.7014
$ws.Range("I132").Value = 2256.5
$ws.Range("K132").Value = 6769.5
$ws.Range("M132").Value = -4239.5
# Row 136
$ws.Range("H136").Value = 4591.7144
$ws.Range("I136").Value = 2491.1428
$ws.Range("K136").Value = 7473.428400000001
$ws.Range("M136").Value = -4923.428400000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 576.5714
$ws.Range("I113").Value = 467.6
$ws.Range("J113").Value = 849
$ws.Range("K113").Value = 1402.8
$ws.Range("L113").Value = 2547
$ws.Range("M113").Value = 767.1999999999998
$ws.Range("N113").Value = -6887
# Row 122
$ws.Range("H122").Value = 2193.875
$ws.Range("I122").Value = 2253.4666
$ws.Range("K122").Value = 6760.399800000001
$ws.Range("M122").Value = -4310.399800000001
# Row 126
$ws.Range("H126").Value = 9281.571
$ws.Range("I126").Value = 8993.5
$ws.Range("J126").Value = 9396.799999999999
$ws.Range("K126").Value = 26980.5
$ws.Range("L126").Value = 28190.4
$ws.Range("M126").Value = -24510.5
$ws.Range("N126").Value = -33130.39999999999
# Row 128
$ws.Range("H128").Value = 165852.25
$ws.Range("J128").Value = 165852.25
$ws.Range("L128").Value = 165852.25
$ws.Range("N128").Value = -175812.25
# Row 130
$ws.Range("H130").Value = 122979
$ws.Range("J130").Value = 122979
$ws.Range("L130").Value = 122979
$ws.Range("N130").Value = -133019
# Row 132
$ws.Range("H132").Value = 5233.4414
$ws.Range("I132").Value = 4620.25
$ws.Range("J132").Value = 6705.1
$ws.Range("K132").Value = 13860.75
$ws.Range("L132").Value = 20115.3
$ws.Range("M132").Value = -11330.75
$ws.Range("N132").Value = -25175.3
# Row 136
$ws.Range("H136").Value = 29466.352
$ws.Range("I136").Value = 1874.52
$ws.Range("K136").Value = 5623.559999999999
$ws.Range("M136").Value = -3073.559999999999

